# Updates cryptos list figures (Price column D, Volume(1h) column E)
# for Sheet1, matching the latest scrape snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.954.71"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.674.86"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'215.21"
$ws.Range("D6").Value = "'0.517"
$ws.Range("E6").Value = "  +1.83%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "1.911.14"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "1.684.43"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").Value = "'65.64"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "26.970.64"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "'8.11"
$ws.Range("E18").Value = "  +4.33%  "
$ws.Range("D19").Value = "'235.30"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'4.45"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'9.17"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("D25").Value = "'145.70"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").Value = "'7.19"
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("D27").Value = "'16.06"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("D33").Value = "1.486.22"
$ws.Range("E33").Value = "  -4.74%  "
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "'0.584"
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("E40").Value = "  +8.15%  "
$ws.Range("D41").Value = "'5.85"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("D45").Value = "1.819.29"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").Value = "'0.779"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").Value = "'90.71"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("E51").Value = "  +0.32%  "
